$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue ($ws.Range("D2")) '27.213.05'
Set-TextValue ($ws.Range("E2")) '  -0.55%  '
Set-TextValue ($ws.Range("D3")) '1.646.70'
Set-TextValue ($ws.Range("E3")) '  -0.94%  '
Set-TextValue ($ws.Range("E4")) '  -0.20%  '
Set-TextValue ($ws.Range("D5")) '217.40'
Set-TextValue ($ws.Range("E5")) '  -1.27%  '
Set-TextValue ($ws.Range("D6")) '0.511'
Set-TextValue ($ws.Range("E6")) '  +0.45%  '
Set-TextValue ($ws.Range("E7")) '  -0.19%  '
Set-TextValue ($ws.Range("E8")) '  +0.16%  '
Set-TextValue ($ws.Range("E9")) '  +0.30%  '
Set-TextValue ($ws.Range("D10")) '19.97'
Set-TextValue ($ws.Range("E10")) '  -0.21%  '
Set-TextValue ($ws.Range("E11")) '  -0.66%  '
Set-TextValue ($ws.Range("D12")) '1.876.71'
Set-TextValue ($ws.Range("E12")) '  -0.96%  '
Set-TextValue ($ws.Range("D13")) '1.654.41'
Set-TextValue ($ws.Range("E13")) '  -0.61%  '
Set-TextValue ($ws.Range("E14")) '  -2.41%  '
Set-TextValue ($ws.Range("D15")) '0.540'
Set-TextValue ($ws.Range("D16")) '67.43'
Set-TextValue ($ws.Range("E16")) '  +0.07%  '
Set-TextValue ($ws.Range("D17")) '27.197.21'
Set-TextValue ($ws.Range("E17")) '  -0.53%  '
Set-TextValue ($ws.Range("D18")) '0.0₃0740'
Set-TextValue ($ws.Range("E18")) '  +0.30%  '
Set-TextValue ($ws.Range("D19")) '218.79'
Set-TextValue ($ws.Range("E19")) '  -1.78%  '
Set-TextValue ($ws.Range("E20")) '  -0.22%  '
Set-TextValue ($ws.Range("D21")) '6.82'
Set-TextValue ($ws.Range("E21")) '  +0.22%  '
Set-TextValue ($ws.Range("D22")) '4.45'
Set-TextValue ($ws.Range("E22")) '  -0.29%  '
Set-TextValue ($ws.Range("E23")) '  -0.08%  '
Set-TextValue ($ws.Range("D24")) '9.19'
Set-TextValue ($ws.Range("E24")) '  -1.25%  '
Set-TextValue ($ws.Range("D25")) '147.50'
Set-TextValue ($ws.Range("E25")) '  +0.14%  '
Set-TextValue ($ws.Range("E26")) '  -0.21%  '
Set-TextValue ($ws.Range("D27")) '7.49'
Set-TextValue ($ws.Range("E27")) '  +0.46%  '
Set-TextValue ($ws.Range("E28")) '  -1.15%  '
Set-TextValue ($ws.Range("D29")) '15.78'
Set-TextValue ($ws.Range("E29")) '  -1.83%  '
Set-TextValue ($ws.Range("D30")) '0.0506'
Set-TextValue ($ws.Range("E30")) '  -1.89%  '
Set-TextValue ($ws.Range("E31")) '  -1.36%  '
Set-TextValue ($ws.Range("E32")) '  -0.85%  '
Set-TextValue ($ws.Range("E33")) '  +0.51%  '
Set-TextValue ($ws.Range("D34")) '1.59'
Set-TextValue ($ws.Range("E34")) '  +1.27%  '
Set-TextValue ($ws.Range("D35")) '1.269.25'
Set-TextValue ($ws.Range("E35")) '  +0.25%  '
Set-TextValue ($ws.Range("E36")) '  +0.12%  '
Set-TextValue ($ws.Range("D37")) '0.0178'
Set-TextValue ($ws.Range("E37")) '  +0.25%  '
Set-TextValue ($ws.Range("D38")) '0.543'
Set-TextValue ($ws.Range("E38")) '  +0.94%  '
Set-TextValue ($ws.Range("E39")) '  +0.47%  '
Set-TextValue ($ws.Range("E40")) '  -0.25%  '
Set-TextValue ($ws.Range("D41")) '0.810'
Set-TextValue ($ws.Range("E41")) '  -0.76%  '
Set-TextValue ($ws.Range("E42")) '  +4.13%  '
Set-TextValue ($ws.Range("D43")) '5.42'
Set-TextValue ($ws.Range("E43")) '  +0.25%  '
Set-TextValue ($ws.Range("D44")) '1.787.05'
Set-TextValue ($ws.Range("E44")) '  -1.09%  '
Set-TextValue ($ws.Range("D45")) '62.43'
Set-TextValue ($ws.Range("E45")) '  +0.91%  '
Set-TextValue ($ws.Range("D46")) '91.80'
Set-TextValue ($ws.Range("E46")) '  -0.89%  '
Set-TextValue ($ws.Range("E47")) '  -0.82%  '
Set-TextValue ($ws.Range("D48")) '0.0₆0106'
Set-TextValue ($ws.Range("E48")) '  +16.51%  '
Set-TextValue ($ws.Range("E49")) '  -1.64%  '
Set-TextValue ($ws.Range("D50")) '7.71'
Set-TextValue ($ws.Range("E50")) '  +0.29%  '
Set-TextValue ($ws.Range("D51")) '0.0973'
Set-TextValue ($ws.Range("E51")) '  -1.32%  '
